# Added a new "Linear Algebra" course row to the Courses sheet
# (student/faculty view now also shows this new offering).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Courses")

$row = 12
$ws.Cells.Item($row, 1).Value = 12          # Course Code
$ws.Cells.Item($row, 2).Value = "Linear Algebra"   # Course Name
$ws.Cells.Item($row, 3).Value = "MATH001"   # Subject Code
$ws.Cells.Item($row, 4).Value = ""          # Section ID (not yet assigned)
$ws.Cells.Item($row, 5).Value = 30          # Capacity
$ws.Cells.Item($row, 6).Value = ""          # Meeting Days/Time
$ws.Cells.Item($row, 7).Value = ""          # Final Exam Date/Time
$ws.Cells.Item($row, 8).Value = ""          # Location
$ws.Cells.Item($row, 9).Value = "Dr. Alan Turing"   # Instructor
